$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.648.21"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.503.24"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.04"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.86"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "3.500.70"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("E11").Value = "  +6.88%  "
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "4.095.57"
$ws.Range("D16").Value = "3.514.72"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "67.561.37"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.52"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.87"
$ws.Range("E21").Value = "  +4.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.83"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.23"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").Value = "3.642.97"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  -1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.75"
$ws.Range("E29").Value = "  +4.87%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  +6.63%  "
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.996"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.15"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "3.493.13"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("E40").Value = "  +7.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.99"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.46"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.892"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "30.15"
$ws.Range("E46").Value = "  +5.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.36"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +2.48%  "
